# Update "last generated" timestamps in the handback-status report.
# Sheet "Overview": Correspond Handback DateTime for the
#   dde28d1d-... (.md) row moves from 06:43:12 -> 06:44:11
# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime
#   for the dde28d1d-... row move from 06:43:07 -> 06:43:58
#   and 06:43:35 -> 06:44:31
# Sheet "de-de": Correspond Handoff Datetime for the
#   dde28d1d-... row moves from 06:43:42 -> 06:44:38

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 06:44:11"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-18 06:43:58"
$wsZhCn.Range("K3").Value = "2016-08-18 06:44:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-18 06:44:38"
